# Insert a new data row at row 261 (pushing existing rows 261..303 down to 262..304)
# and populate it with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 261, shifting rows 261-303 down to 262-304.
$ws.Rows.Item(261).Insert()

# Populate the newly inserted row 261 with the new record's values.
$ws.Cells.Item(261, 1).Value = 8
$ws.Cells.Item(261, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(261, 3).Value = "Coquimbo"
$ws.Cells.Item(261, 4).Value = 45180
$ws.Cells.Item(261, 5).Value = 4
$ws.Cells.Item(261, 6).Value = 100112001
$ws.Cells.Item(261, 7).Value = "Berenjena"
$ws.Cells.Item(261, 8).Value = "Sin especificar"
$ws.Cells.Item(261, 9).Value = "Primera"
$ws.Cells.Item(261, 10).Value = 460
$ws.Cells.Item(261, 11).Value = 9000
$ws.Cells.Item(261, 12).Value = 10000
$ws.Cells.Item(261, 13).Value = 9500
$ws.Cells.Item(261, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(261, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(261, 16).Value = 190
$ws.Cells.Item(261, 17).Value = 50
$ws.Cells.Item(261, 18).Value = "Hortaliza"
